$d = $word.ActiveDocument

# 1) Locate the paragraph containing "Name: title}}" and strip the
#    stray "}}" left over from the broken template-key handling.
$p1 = $d.Paragraphs(1)
$r1 = $p1.Range
$r1.Find.Execute("Name: title}}", $true, $false, $false, $false, $false, `
                 $true, 1, $false, "Name: title", 2)

# 2) Insert a new paragraph right after it reporting the broken key,
#    matching the pre-"fix" behaviour being reverted.
$p1 = $d.Paragraphs(1)
$p1.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs(2)
$p2.Range.InsertBefore("Broken template key: title")
